$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sector/correlation data (rows 2-39 after edit)
$data = @(
    @('Multi-Utilities(18)', 0.6117445464417549),
    @('Energy Equipment & Services(32)', 0.502640014568374),
    @('Gas Utilities(12)', 0.4821385972405633),
    @('Electric Utilities(28)', 0.4674867790683934),
    @('Containers & Packaging(12)', 0.456085570281568),
    @('Construction & Engineering(20)', 0.4228573107356991),
    @('Auto Components(21)', 0.4138606585033447),
    @('Road & Rail(22)', 0.3910516778365106),
    @('Water Utilities(12)', 0.3904295022997293),
    @('Airlines(14)', 0.3735498953362092),
    @('Insurance(75)', 0.3684458161108561),
    @('Marine(15)', 0.3549270385176713),
    @('Building Products(23)', 0.3459578419204664),
    @('Machinery(85)', 0.3177504155016216),
    @('Consumer Finance(15)', 0.313964282229209),
    @('Trading Companies & Distributors(25)', 0.2858472029886769),
    @('Oil, Gas & Consumable Fuels(122)', 0.2646615236424472),
    @('Metals & Mining(89)', 0.2646001487063055),
    @('Specialty Retail(58)', 0.2635153245892091),
    @('Banks(246)', 0.2521303549963045),
    @('Beverages(21)', 0.2331922912862205),
    @('Diversified Telecommunication Services(20)', 0.2212130460900358),
    @('Food Products(44)', 0.2174162042938039),
    @('Household Durables(39)', 0.2172100384495448),
    @('Commercial Services & Supplies(52)', 0.2150374974613652),
    @('Media(42)', 0.2050837148192892),
    @('Health Care Providers & Services(46)', 0.19985448571264),
    @('Semiconductors & Semiconductor Equipment(68)', 0.1868380093077185),
    @('Electrical Equipment(28)', 0.1786006896993373),
    @('Health Care Equipment & Supplies(83)', 0.1747394250851111),
    @('Chemicals(51)', 0.1671865513529269),
    @('Capital Markets(75)', 0.1610382427600499),
    @('Thrifts & Mortgage Finance(47)', 0.1543685246506317),
    @('Hotels, Restaurants & Leisure(50)', 0.1439368399565081),
    @('Software(66)', 0.1436270686000481),
    @('IT Services(52)', 0.1288975452186246),
    @('Professional Services(35)', 0.1206784113690283),
    @('Biotechnology(126)', 0.1169280115719262)
)

$oldLastRow = 42
$newLastRow = 1 + $data.Count   # 39

# Write the new values into rows 2..39
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-obsolete trailing rows (old rows 40-42) since the table shrank
if ($newLastRow -lt $oldLastRow) {
    $clearRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 2))
    $clearRange.ClearContents()
}

$wb.Save()
